# update scripts wuth new tpm
#
# NATMI Nodal-Tdgf1 ligand-receptor pair sheet was regenerated from new TPM
# values. This:
#   - recomputes the NATMI statistics for the existing Sending-cluster rows
#     (ECs / FAPs / MuSCs->Inflammatory-Mac) against Nodal/Tdgf1/ECs,
#   - changes the Target cluster from "Resolving-Mac" to "ECs" for every row,
#   - adds a new trailing row for the "MuSCs" sending cluster (previously the
#     last row, now pushed down since the macrophage-state label changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sending cluster: ECs) ---------------------------------------
$ws.Range("B2").Value = "Nodal"
$ws.Range("C2").Value = "Tdgf1"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.7356009999999999
$ws.Range("H2").Value = 2.206803
$ws.Range("I2").Value = 0.5568025227146887
$ws.Range("J2").Value = 0.5688250870390384
$ws.Range("M2").Value = 0.08586166666666667
$ws.Range("N2").Value = 0.257585
$ws.Range("Q2").Value = 0.06315992786166666
$ws.Range("R2").Value = 0.568439350755
$ws.Range("S2").Value = 0.5568025227146887
$ws.Range("T2").Value = 0.5688250870390384

# --- Row 3 (Sending cluster: FAPs) --------------------------------------
$ws.Range("B3").Value = "Nodal"
$ws.Range("C3").Value = "Tdgf1"
$ws.Range("D3").Value = "ECs"
$ws.Range("H3").Value = 0.9457549999999999
$ws.Range("I3").Value = 0.2386251830680085
$ws.Range("J3").Value = 0.2437776141289484
$ws.Range("M3").Value = 0.08586166666666667
$ws.Range("N3").Value = 0.257585
$ws.Range("Q3").Value = 0.02706803351944444
$ws.Range("R3").Value = 0.243612301675
$ws.Range("S3").Value = 0.2386251830680085
$ws.Range("T3").Value = 0.2437776141289484

# --- Row 4 (Sending cluster becomes: Inflammatory-Mac) -------------------
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Nodal"
$ws.Range("C4").Value = "Tdgf1"
$ws.Range("D4").Value = "ECs"
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1864953333333333
$ws.Range("H4").Value = 0.559486
$ws.Range("I4").Value = 0.1411649414213912
$ws.Range("J4").Value = 0.1442129961972698
$ws.Range("M4").Value = 0.08586166666666667
$ws.Range("N4").Value = 0.257585
$ws.Range("Q4").Value = 0.01601280014555556
$ws.Range("R4").Value = 0.14411520131
$ws.Range("S4").Value = 0.1411649414213912
$ws.Range("T4").Value = 0.1442129961972698

# --- Row 5 (new row, Sending cluster: MuSCs) ------------------------------
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Nodal"
$ws.Range("C5").Value = "Tdgf1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.0837685
$ws.Range("H5").Value = 0.167537
$ws.Range("I5").Value = 0.06340735279591164
$ws.Range("J5").Value = 0.04318430263474329
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.08586166666666667
$ws.Range("N5").Value = 0.257585
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.007192503024166666
$ws.Range("R5").Value = 0.043155018145
$ws.Range("S5").Value = 0.06340735279591164
$ws.Range("T5").Value = 0.04318430263474329

Write-Output "Updated Nodal-Tdgf1 sheet with new TPM-derived values (rows 2-5)."
